# Realestate Update resale numbers 2024-01-17 09:35
# Appends a new data row (row 66) to the CityResaleNum sheet with the
# latest resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column: force text storage so "2024-01-17" isn't coerced into a
# date serial number (matches the existing text-stored date cells above it).
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "2024-01-17"
$ws.Range("A66").ClearFormats()

# --- Time / Weekday: plain strings, not number-like, stored as text as-is.
$ws.Range("B66").Value = "09:34:59"
$ws.Range("C66").Value = "Wednesday"

# --- Week column: force text storage so "02" keeps its leading zero
# instead of being coerced to the number 2.
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "02"
$ws.Range("D66").ClearFormats()

# --- City resale-number columns (plain numeric values).
$ws.Range("E66").Value = 138780
$ws.Range("F66").Value = 139464
$ws.Range("G66").Value = 170610
$ws.Range("H66").Value = 148418
$ws.Range("I66").Value = -1
$ws.Range("J66").Value = 118909
$ws.Range("K66").Value = 221821
$ws.Range("L66").Value = 254550
$ws.Range("M66").Value = 184944
$ws.Range("N66").Value = 110343
$ws.Range("O66").Value = 41270
$ws.Range("P66").Value = 30924
$ws.Range("Q66").Value = 73337
$ws.Range("R66").Value = -1
$ws.Range("S66").Value = 42146
$ws.Range("T66").Value = -1
